$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "FY1"
$ws.Range("C2").Value = -0.99476
$ws.Range("D2").Value = 0.012181
$ws.Range("E2").Value = -0.101506
$ws.Range("F2").Value = 140
$ws.Range("G2").Value = 8000
$ws.Range("H2").Value = 120
$ws.Range("I2").Value = 800
$ws.Range("J2").Value = 70.368708
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 49.258096
$ws.Range("M2").Value = 10940
$ws.Range("N2").Value = 84
$ws.Range("O2").Value = 1100
$ws.Range("P2").Value = 8000
$ws.Range("Q2").Value = 120
$ws.Range("R2").Value = 80
$ws.Range("S2").Value = 51.258096
$ws.Range("T2").Value = 1.2

$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "FY1"
$ws.Range("C3").Value = -0.99476
$ws.Range("D3").Value = 0.012181
$ws.Range("E3").Value = -0.101506
$ws.Range("F3").Value = 140
$ws.Range("G3").Value = 8000
$ws.Range("H3").Value = 120
$ws.Range("I3").Value = 800
$ws.Range("J3").Value = 70.368708
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 56.294967
$ws.Range("M3").Value = 9960
$ws.Range("N3").Value = 96
$ws.Range("O3").Value = 1000
$ws.Range("P3").Value = 8000
$ws.Range("Q3").Value = 120
$ws.Range("R3").Value = 80
$ws.Range("S3").Value = 58.294967
$ws.Range("T3").Value = 1.26

$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "FY1"
$ws.Range("C4").Value = -0.99476
$ws.Range("D4").Value = 0.012181
$ws.Range("E4").Value = -0.101506
$ws.Range("F4").Value = 140
$ws.Range("G4").Value = 8000
$ws.Range("H4").Value = 120
$ws.Range("I4").Value = 800
$ws.Range("J4").Value = 70.368708
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 63.331837
$ws.Range("M4").Value = 8980
$ws.Range("N4").Value = 108
$ws.Range("O4").Value = 900
$ws.Range("P4").Value = 8000
$ws.Range("Q4").Value = 120
$ws.Range("R4").Value = 80
$ws.Range("S4").Value = 65.33183699999999
$ws.Range("T4").Value = 1.32

$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "FY2"
$ws.Range("C5").Value = -0.960422
$ws.Range("D5").Value = -0.243947
$ws.Range("E5").Value = -0.134459
$ws.Range("F5").Value = 140
$ws.Range("G5").Value = 7000
$ws.Range("H5").Value = 130
$ws.Range("I5").Value = 700
$ws.Range("J5").Value = 37.186033
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 26.030223
$ws.Range("M5").Value = 8500
$ws.Range("N5").Value = 511
$ws.Range("O5").Value = 910
$ws.Range("P5").Value = 7000
$ws.Range("Q5").Value = 130
$ws.Range("R5").Value = 80
$ws.Range("S5").Value = 28.030223
$ws.Range("T5").Value = 1.2

$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "FY2"
$ws.Range("C6").Value = -0.960422
$ws.Range("D6").Value = -0.243947
$ws.Range("E6").Value = -0.134459
$ws.Range("F6").Value = 140
$ws.Range("G6").Value = 7000
$ws.Range("H6").Value = 130
$ws.Range("I6").Value = 700
$ws.Range("J6").Value = 37.186033
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 29.748826
$ws.Range("M6").Value = 8000
$ws.Range("N6").Value = 384
$ws.Range("O6").Value = 840
$ws.Range("P6").Value = 7000
$ws.Range("Q6").Value = 130
$ws.Range("R6").Value = 80
$ws.Range("S6").Value = 31.748826
$ws.Range("T6").Value = 1.26

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "FY3"
$ws.Range("C7").Value = -0.089975
$ws.Range("D7").Value = 0.995723
$ws.Range("E7").Value = -0.020994
$ws.Range("F7").Value = 87.193557
$ws.Range("G7").Value = 5700
$ws.Range("H7").Value = 320
$ws.Range("I7").Value = 630
$ws.Range("J7").Value = 38.23977
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 26.767839
$ws.Range("M7").Value = 5790
$ws.Range("N7").Value = -676
$ws.Range("O7").Value = 651
$ws.Range("P7").Value = 5700
$ws.Range("Q7").Value = 320
$ws.Range("R7").Value = 80
$ws.Range("S7").Value = 28.767839
$ws.Range("T7").Value = 1.2

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "FY3"
$ws.Range("C8").Value = -0.089975
$ws.Range("D8").Value = 0.995723
$ws.Range("E8").Value = -0.020994
$ws.Range("F8").Value = 87.193557
$ws.Range("G8").Value = 5700
$ws.Range("H8").Value = 320
$ws.Range("I8").Value = 630
$ws.Range("J8").Value = 38.23977
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 30.591816
$ws.Range("M8").Value = 5760
$ws.Range("N8").Value = -344
$ws.Range("O8").Value = 644
$ws.Range("P8").Value = 5700
$ws.Range("Q8").Value = 320
$ws.Range("R8").Value = 80
$ws.Range("S8").Value = 32.591816
$ws.Range("T8").Value = 1.26

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "FY3"
$ws.Range("C9").Value = -0.089975
$ws.Range("D9").Value = 0.995723
$ws.Range("E9").Value = -0.020994
$ws.Range("F9").Value = 87.193557
$ws.Range("G9").Value = 5700
$ws.Range("H9").Value = 320
$ws.Range("I9").Value = 630
$ws.Range("J9").Value = 38.23977
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 34.415793
$ws.Range("M9").Value = 5730
$ws.Range("N9").Value = -12
$ws.Range("O9").Value = 637
$ws.Range("P9").Value = 5700
$ws.Range("Q9").Value = 320
$ws.Range("R9").Value = 80
$ws.Range("S9").Value = 36.415793
$ws.Range("T9").Value = 1.32

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "FY4"
$ws.Range("C10").Value = -0.946773
$ws.Range("D10").Value = -0.257522
$ws.Range("E10").Value = -0.193142
$ws.Range("F10").Value = 140
$ws.Range("G10").Value = 4750
$ws.Range("H10").Value = 300
$ws.Range("I10").Value = 525
$ws.Range("J10").Value = 47.152664
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 33.006865
$ws.Range("M10").Value = 6625
$ws.Range("N10").Value = 810
$ws.Range("O10").Value = 907.5
$ws.Range("P10").Value = 4750
$ws.Range("Q10").Value = 300
$ws.Range("R10").Value = 80
$ws.Range("S10").Value = 35.006865
$ws.Range("T10").Value = 1.2

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "FY4"
$ws.Range("C11").Value = -0.946773
$ws.Range("D11").Value = -0.257522
$ws.Range("E11").Value = -0.193142
$ws.Range("F11").Value = 140
$ws.Range("G11").Value = 4750
$ws.Range("H11").Value = 300
$ws.Range("I11").Value = 525
$ws.Range("J11").Value = 47.152664
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 37.722131
$ws.Range("M11").Value = 6000
$ws.Range("N11").Value = 640
$ws.Range("O11").Value = 780
$ws.Range("P11").Value = 4750
$ws.Range("Q11").Value = 300
$ws.Range("R11").Value = 80
$ws.Range("S11").Value = 39.722131
$ws.Range("T11").Value = 1.26

$ws.Range("A12").Value = "M3"
$ws.Range("B12").Value = "FY5"
$ws.Range("C12").Value = -0.972812
$ws.Range("D12").Value = 0.211121
$ws.Range("E12").Value = -0.095211
$ws.Range("F12").Value = 140
$ws.Range("G12").Value = 3600
$ws.Range("H12").Value = 40
$ws.Range("I12").Value = 380
$ws.Range("J12").Value = 69.01937
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 48.313559
$ws.Range("M12").Value = 6420
$ws.Range("N12").Value = -572
$ws.Range("O12").Value = 656
$ws.Range("P12").Value = 3600
$ws.Range("Q12").Value = 40
$ws.Range("R12").Value = 80
$ws.Range("S12").Value = 50.313559
$ws.Range("T12").Value = 1.2

$ws.Range("A13").Value = "M3"
$ws.Range("B13").Value = "FY5"
$ws.Range("C13").Value = -0.972812
$ws.Range("D13").Value = 0.211121
$ws.Range("E13").Value = -0.095211
$ws.Range("F13").Value = 140
$ws.Range("G13").Value = 3600
$ws.Range("H13").Value = 40
$ws.Range("I13").Value = 380
$ws.Range("J13").Value = 69.01937
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 55.215496
$ws.Range("M13").Value = 5480
$ws.Range("N13").Value = -368
$ws.Range("O13").Value = 564
$ws.Range("P13").Value = 3600
$ws.Range("Q13").Value = 40
$ws.Range("R13").Value = 80
$ws.Range("S13").Value = 57.215496
$ws.Range("T13").Value = 1.26

$ws.Range("A14").Value = "M3"
$ws.Range("B14").Value = "FY5"
$ws.Range("C14").Value = -0.972812
$ws.Range("D14").Value = 0.211121
$ws.Range("E14").Value = -0.095211
$ws.Range("F14").Value = 140
$ws.Range("G14").Value = 3600
$ws.Range("H14").Value = 40
$ws.Range("I14").Value = 380
$ws.Range("J14").Value = 69.01937
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 62.117433
$ws.Range("M14").Value = 4540
$ws.Range("N14").Value = -164
$ws.Range("O14").Value = 472
$ws.Range("P14").Value = 3600
$ws.Range("Q14").Value = 40
$ws.Range("R14").Value = 80
$ws.Range("S14").Value = 64.11743300000001
$ws.Range("T14").Value = 1.32

Write-Host "Updated rows 2-14"